# Auto-applying scheduled market-data refresh to Sheets/Midgardsormr_Profits
# (chore: update Sheets via scheduled runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2965.7273
$ws.Range("I98").Value = 616.4483
$ws.Range("J98").Value = 19998
$ws.Range("K98").Value = 616.4483
$ws.Range("L98").Value = 19998
$ws.Range("M98").Value = 881.5517
$ws.Range("N98").Value = -22994
$ws.Range("H112").Value = 6625.456
$ws.Range("I112").Value = 1666.5
$ws.Range("J112").Value = 6805.7817
$ws.Range("K112").Value = 4999.5
$ws.Range("L112").Value = 20417.3451
$ws.Range("M112").Value = -3891.5
$ws.Range("N112").Value = -22633.3451
$ws.Range("H122").Value = 2965.7273
$ws.Range("I122").Value = 616.4483
$ws.Range("J122").Value = 19998
$ws.Range("K122").Value = 1849.3449
$ws.Range("L122").Value = 59994
$ws.Range("M122").Value = 600.6550999999999
$ws.Range("N122").Value = -64894
$ws.Range("H137").Value = 7546.403
$ws.Range("I137").Value = 10311.211
$ws.Range("J137").Value = 3923.5518
$ws.Range("K137").Value = 30933.633
$ws.Range("L137").Value = 11770.6554
$ws.Range("M137").Value = -28383.633
$ws.Range("N137").Value = -16870.6554

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6676.8687
$ws.Range("I32").Value = 6603.3613
$ws.Range("K32").Value = 6603.3613
$ws.Range("M32").Value = -6316.3613
$ws.Range("H45").Value = 4832.76
$ws.Range("I45").Value = 3364.6667
$ws.Range("J45").Value = 6187.923
$ws.Range("K45").Value = 3364.6667
$ws.Range("L45").Value = 6187.923
$ws.Range("M45").Value = -2987.6667
$ws.Range("N45").Value = -6941.923
$ws.Range("H74").Value = 143971.4
$ws.Range("I74").Value = 167809.47
$ws.Range("J74").Value = 943
$ws.Range("K74").Value = 167809.47
$ws.Range("L74").Value = 943
$ws.Range("M74").Value = -166935.47
$ws.Range("N74").Value = -2691
$ws.Range("H77").Value = 143971.4
$ws.Range("I77").Value = 167809.47
$ws.Range("J77").Value = 943
$ws.Range("K77").Value = 839047.35
$ws.Range("L77").Value = 4715
$ws.Range("M77").Value = -834679.35
$ws.Range("N77").Value = -13451
$ws.Range("H110").Value = 1392.625
$ws.Range("I110").Value = 1504.6666
$ws.Range("J110").Value = 1056.5
$ws.Range("K110").Value = 1504.6666
$ws.Range("L110").Value = 1056.5
$ws.Range("M110").Value = 540.3334
$ws.Range("N110").Value = -5146.5
$ws.Range("H132").Value = 1236.5167
$ws.Range("I132").Value = 1185.569
$ws.Range("J132").Value = 2714
$ws.Range("K132").Value = 3556.707
$ws.Range("L132").Value = 8142
$ws.Range("M132").Value = -1026.707
$ws.Range("N132").Value = -13202
$ws.Range("H134").Value = 58126.332
$ws.Range("J134").Value = 58126.332
$ws.Range("L134").Value = 58126.332
$ws.Range("N134").Value = -68266.33199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 20420.111
$ws.Range("I20").Value = 31212.883
$ws.Range("J20").Value = 2072.4
$ws.Range("K20").Value = 31212.883
$ws.Range("L20").Value = 2072.4
$ws.Range("M20").Value = -30965.883
$ws.Range("N20").Value = -2566.4
$ws.Range("H138").Value = 147512
$ws.Range("J138").Value = 147512
$ws.Range("L138").Value = 147512
$ws.Range("N138").Value = -157792

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 135.24
$ws.Range("I7").Value = 149.64706
$ws.Range("K7").Value = 149.64706
$ws.Range("M7").Value = -36.64706000000001
$ws.Range("H22").Value = 869.36365
$ws.Range("I22").Value = 888.1667
$ws.Range("J22").Value = 846.8
$ws.Range("K22").Value = 888.1667
$ws.Range("L22").Value = 846.8
$ws.Range("M22").Value = -538.1667
$ws.Range("N22").Value = -1546.8
$ws.Range("H86").Value = 28147.133
$ws.Range("I86").Value = 40244.312
$ws.Range("J86").Value = 14321.786
$ws.Range("K86").Value = 40244.312
$ws.Range("L86").Value = 14321.786
$ws.Range("M86").Value = -39121.312
$ws.Range("N86").Value = -16567.786
$ws.Range("H89").Value = 28147.133
$ws.Range("I89").Value = 40244.312
$ws.Range("J89").Value = 14321.786
$ws.Range("K89").Value = 201221.56
$ws.Range("L89").Value = 71608.92999999999
$ws.Range("M89").Value = -195605.56
$ws.Range("N89").Value = -82840.92999999999
$ws.Range("H105").Value = 1397.1428
$ws.Range("I105").Value = 1499.5834
$ws.Range("J105").Value = 782.5
$ws.Range("K105").Value = 1499.5834
$ws.Range("L105").Value = 782.5
$ws.Range("M105").Value = 247.4166
$ws.Range("N105").Value = -4276.5
$ws.Range("H132").Value = 13961.872
$ws.Range("I132").Value = 13961.872
$ws.Range("K132").Value = 41885.61599999999
$ws.Range("M132").Value = -39355.61599999999
$ws.Range("H141").Value = 52312.3
$ws.Range("J141").Value = 52312.3
$ws.Range("L141").Value = 52312.3
$ws.Range("N141").Value = -62672.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 248
$ws.Range("J12").Value = 221.15
$ws.Range("L12").Value = 663.45
$ws.Range("N12").Value = -1009.45
$ws.Range("H23").Value = 293
$ws.Range("I23").Value = 4
$ws.Range("J23").Value = 389.33334
$ws.Range("K23").Value = 12
$ws.Range("L23").Value = 1168.00002
$ws.Range("M23").Value = 223
$ws.Range("N23").Value = -1638.00002
$ws.Range("H68").Value = 4365.433
$ws.Range("I68").Value = 1520
$ws.Range("J68").Value = 4681.593
$ws.Range("K68").Value = 4560
$ws.Range("L68").Value = 14044.779
$ws.Range("M68").Value = -3749
$ws.Range("N68").Value = -15666.779
$ws.Range("H71").Value = 4365.433
$ws.Range("I71").Value = 1520
$ws.Range("J71").Value = 4681.593
$ws.Range("K71").Value = 13680
$ws.Range("L71").Value = 42134.337
$ws.Range("M71").Value = -9624
$ws.Range("N71").Value = -50246.337
$ws.Range("H81").Value = 3778.3333
$ws.Range("I81").Value = 1667.5
$ws.Range("K81").Value = 5002.5
$ws.Range("M81").Value = -3879.5
$ws.Range("H84").Value = 3778.3333
$ws.Range("I84").Value = 1667.5
$ws.Range("K84").Value = 15007.5
$ws.Range("M84").Value = -9391.5
$ws.Range("H129").Value = 2640.7932
$ws.Range("I129").Value = 1298.6666
$ws.Range("J129").Value = 3588.1765
$ws.Range("K129").Value = 3895.9998
$ws.Range("L129").Value = 10764.5295
$ws.Range("M129").Value = 1104.0002
$ws.Range("N129").Value = -20764.5295
$ws.Range("H140").Value = 2930.077
$ws.Range("I140").Value = 2930.077
$ws.Range("K140").Value = 8790.231
$ws.Range("M140").Value = -3610.231
$ws.Range("H141").Value = 4708.522
$ws.Range("I141").Value = 4635.4116
$ws.Range("K141").Value = 13906.2348
$ws.Range("M141").Value = -8726.234800000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7400
$ws.Range("I70").Value = 7000
$ws.Range("K70").Value = 7000
$ws.Range("M70").Value = -6730
$ws.Range("H73").Value = 7400
$ws.Range("I73").Value = 7000
$ws.Range("K73").Value = 7000
$ws.Range("M73").Value = -6064
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("H136").Value = 9038.031999999999
$ws.Range("J136").Value = 9038.031999999999
$ws.Range("L136").Value = 27114.096
$ws.Range("N136").Value = -32214.096
$ws.Range("H138").Value = 52586.668
$ws.Range("J138").Value = 52586.668
$ws.Range("L138").Value = 52586.668
$ws.Range("N138").Value = -62866.668
$ws.Range("N87").ClearContents()
$ws.Range("N90").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 8033.2
$ws.Range("I46").Value = 3483.6667
$ws.Range("J46").Value = 9983
$ws.Range("K46").Value = 3483.6667
$ws.Range("L46").Value = 9983
$ws.Range("M46").Value = -3295.6667
$ws.Range("N46").Value = -10359

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 65000
$ws.Range("J86").Value = 65000
$ws.Range("L86").Value = 65000
$ws.Range("N86").Value = -67246
$ws.Range("H89").Value = 65000
$ws.Range("J89").Value = 65000
$ws.Range("L89").Value = 325000
$ws.Range("N89").Value = -336232
$ws.Range("H122").Value = 15193442
$ws.Range("I122").Value = 20053910
$ws.Range("K122").Value = 60161730
$ws.Range("M122").Value = -60159280
$ws.Range("H132").Value = 4645788
$ws.Range("I132").Value = 5017051.5
$ws.Range("J132").Value = 4996.75
$ws.Range("K132").Value = 15051154.5
$ws.Range("L132").Value = 14990.25
$ws.Range("M132").Value = -15048624.5
$ws.Range("N132").Value = -20050.25
